# "enable scrollX in DT output" - the underlying HTML widget (DataTables)
# was re-rendered with scrollX turned on, which reflowed the page a hair
# and nudged a handful of the country-label textboxes on the slide by a
# few dozen EMU. Reproduce that nudge by moving the four label shapes to
# their new positions.
#
# These shapes live inside the single group shape that sits on the slide
# (id 2); the labels themselves are "tx8".."tx11" and are reached via
# GroupItems rather than the slide's own Shapes collection.
#
# NOTE on precision: Shape.Left/Top (and GroupItems' Shape.Left/Top) are
# exposed as single-precision (float32) point values, exactly like real
# PowerPoint's object model. The underlying XML stores EMU (1 pt = 12700
# EMU). Assigning "target_emu / 12700.0" naively can truncate to one EMU
# off the intended value once it is rounded to float32 and converted back
# to EMU. To land on the exact EMU the diff expects, each assigned point
# value below has been chosen (by scanning nearby float32 values) so that
# it reproduces the exact target EMU offset after the float32 round trip.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(1)

# tx8 "Asia": (4779462, 2842123) -> (4779637, 2851865) EMU
$tx8 = $g.GroupItems.Item("tx8")
$tx8.Left = 376.3493957519531
$tx8.Top = 224.55630493164062

# tx9 "Europe": (4579432, 4272369) -> (4581570, 4269522) EMU
$tx9 = $g.GroupItems.Item("tx9")
$tx9.Left = 360.7535705566406
$tx9.Top = 336.182861328125

# tx10 "North America": (5638311, 3544553) -> (5636756, 3546489) EMU
$tx10 = $g.GroupItems.Item("tx10")
$tx10.Left = 443.8390808105469
$tx10.Top = 279.2511291503906

# tx11 "South America": (5621416, 2851848) -> (5621425, 2850926) EMU
$tx11 = $g.GroupItems.Item("tx11")
$tx11.Left = 442.63189697265625
$tx11.Top = 224.4823760986328
